# "Add files via upload" - fills in the missing ID / S/N pair for the
# worker in row 13 (Jeffrey Smith), who previously had no C/D values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "3A5DB840"
$ws.Range("D13").Value = "N521D5060006"

# Leave the selection on the newly-filled cell, matching where the user's
# cursor ended up after typing the value in.
$ws.Range("D13").Select() | Out-Null
